# Deploy the implementation guide.
#
# 1. Metadata!B6 ("Status" value): active -> draft
# 2. Metadata!B8 ("Date" value): refreshed publish timestamp
# 3. Re-assert the header/body wrap formatting (vertical=top, wrap text)
#    on every cell of both sheets so the saved styles keep their
#    alignment formatting (and pick up the applyAlignment flag) instead
#    of losing it on save.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

foreach ($ws in @($meta, $concepts)) {
    $used = $ws.UsedRange
    $used.VerticalAlignment = -4160  # xlVAlignTop
    $used.WrapText = $true
}
